# Omaha_Cal_Info_CP02PMCI_00001.xlsx — "Adding updated cal sheets."
#
# Logical changes applied:
#   * Moorings!A2 "Ref Des" value corrected: CP02PMCI -> CP02PMCI-PM001
#   * Asset_Cal_Info!E23:E26 calibration-coefficient names: underscores -> spaces
#     (CC_angular_resolution -> CC_angular resolution, etc.)
#   * Asset_Cal_Info!E34:E35 calibration-coefficient names replaced with the
#     new ones (CC_dark_offset / CC_scale_wet) and highlighted in red, 10pt,
#     left aligned, to flag the update.
#   * Workbook-level calculation option (multi-threaded calc) toggled off.
#   * Selection/view state restored to match the saved file (best effort).

$wb = $excel.ActiveWorkbook

# --- Workbook level calc option (best effort; mirrors concurrentCalc="0") ---
try {
    $excel.MultiThreadedCalculation.Enabled = $false
} catch {
}

# ---------------------------------------------------------------------------
# Sheet "Moorings"
# ---------------------------------------------------------------------------
$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Range("A2").Value = "CP02PMCI-PM001"

# ---------------------------------------------------------------------------
# Sheet "Asset_Cal_Info"
# ---------------------------------------------------------------------------
$cal = $wb.Worksheets.Item("Asset_Cal_Info")

# Rename the optical-backscatter calibration coefficients to the
# space-separated spellings used going forward.
$cal.Range("E23").Value = "CC_angular resolution"
$cal.Range("E24").Value = "CC_depolarization ratio"
$cal.Range("E25").Value = "CC_measurement wavelength"
$cal.Range("E26").Value = "CC_scattering angle"

# Replace the two fluorometer coefficient names with the updated ones and
# flag them in red so the change stands out on the updated cal sheet.
$cal.Range("E34").Value = "CC_dark_offset"
$cal.Range("E35").Value = "CC_scale_wet"
$flagged = $cal.Range("E34:E35")
$flagged.Font.Size = 10
$flagged.Font.Color = 255
$flagged.HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Restore the on-screen selection / active sheet to match the saved state.
# ---------------------------------------------------------------------------
$moorings.Activate()
$moorings.Range("H2").Select()

$cal.Activate()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 28
    $win.ScrollColumn = 1
} catch {
}
$cal.Range("E42").Select()
